$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.605.09"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "2.529.88"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.97"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -1.57%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -3.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.59"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.34"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "2.915.70"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.26"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.38%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.513.08"
$ws.Range("E16").Value = "  -4.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.808"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.54%  "
$ws.Range("D18").Value = "42.592.81"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.69"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.16"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.39"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.57"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.86"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("E27").Value = "  -6.18%  "
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.11"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.22"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.19"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.73"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("E33").Value = "  +9.98%  "
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.99"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.13"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -8.87%  "
$ws.Range("E38").Value = "  -7.94%  "
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.20"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.43"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.37%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "2.007.43"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.90"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "2.770.58"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "79.14"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("E50").Value = "  -4.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.07"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.08%  "
